# Update 2021 HWL2 First Batch
#
# Extends every per-year data table in the workbook with columns for the
# years 2016-2050 (the source only went up to 2015).
#
# - "Data Clio Infra Format" is the wide table: one column per year. The
#   last year column is 2015 (column TD); 35 new year columns (2016-2050)
#   are appended right after it.
# - "Data Long Format" ends with a "year"/"value" column pair. 35 new
#   (placeholder) year columns are inserted just before that pair, mirroring
#   the wide sheet's layout; the former "year"/"value" pair slides right.
# - "Metadata" has no per-year columns and needs no structural change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Data Clio Infra Format"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Data Clio Infra Format")

$firstNewCol = 525   # column TE, right after TD (=2015)
$lastNewCol  = 559   # column UM

$hdr1 = $ws1.Range($ws1.Cells.Item(1, $firstNewCol), $ws1.Cells.Item(1, $lastNewCol))
$hdr1.NumberFormat = "@"
for ($col = $firstNewCol; $col -le $lastNewCol; $col++) {
    $year = 2016 + ($col - $firstNewCol)
    $ws1.Cells.Item(1, $col).Value = "$year"
}
$hdr1.ClearFormats()

# ---------------------------------------------------------------------------
# Sheet: "Data Long Format"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Data Long Format")

# Push the trailing "year"/"value" columns (E:F) to the right, opening up
# 35 blank columns (E:AM) in their place.
$ws2.Columns("E:AM").Insert()

$hdr2 = $ws2.Range($ws2.Cells.Item(1, 5), $ws2.Cells.Item(1, 39))
$hdr2.NumberFormat = "@"
for ($col = 5; $col -le 39; $col++) {
    $year = 2016 + ($col - 5)
    $ws2.Cells.Item(1, $col).Value = "$year"
}
$hdr2.ClearFormats()
